$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Covid-19 podatki")

$newRowNum = 76

# --- Format the new row to match the row above it (row 75) ---
# Column A: date column (custom date format, font "Calibri Light" themed, top-aligned)
$a = $ws.Cells.Item($newRowNum, 1)
$a.NumberFormat = "d/\ m/\ yyyy;@"
$a.Font.Name = "Calibri Light"
$a.Font.Size = 10
$a.Font.ThemeFont = 2
$a.HorizontalAlignment = -4152
$a.VerticalAlignment = -4160

# Column B: thousands-style number column, right aligned
$b = $ws.Cells.Item($newRowNum, 2)
$b.NumberFormat = "#,##0"
$b.Font.Name = "Calibri Light"
$b.Font.Size = 10
$b.Font.ThemeFont = 2
$b.HorizontalAlignment = -4152

# Columns C-J: general number columns, right aligned
for ($col = 3; $col -le 10; $col++) {
    $c = $ws.Cells.Item($newRowNum, $col)
    $c.Font.Name = "Calibri Light"
    $c.Font.Size = 10
    $c.Font.ThemeFont = 2
    $c.HorizontalAlignment = -4152
}

# --- Write the new day's data (2020-05-25, serial 43976) ---
$ws.Cells.Item($newRowNum, 1).Value = 43976
$ws.Cells.Item($newRowNum, 2).Value = 75770
$ws.Cells.Item($newRowNum, 3).Value = 754
$ws.Cells.Item($newRowNum, 4).Value = 1469
$ws.Cells.Item($newRowNum, 5).Value = 0
$ws.Cells.Item($newRowNum, 6).Value = 9
$ws.Cells.Item($newRowNum, 7).Value = 2
$ws.Cells.Item($newRowNum, 8).Value = 6
$ws.Cells.Item($newRowNum, 9).Value = 108
$ws.Cells.Item($newRowNum, 10).Value = 1

# --- Grow the table (ListObject) so it covers the new row ---
$tbl = $ws.ListObjects.Item("Tabela1")
$tbl.Resize($ws.Range("A1:J76"))

# --- Update the selection to the new row, matching the edit that added it ---
$ws.Range("A76:J76").Select()

$wb.Save()
